# This workbook contains a weekly price log for "Cilantro" at the
# "Macroferia Regional de Talca" market. A new observation was recorded
# for the week, which pushes all the later rows down by one, and a
# further new observation was appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; this shifts the existing rows
# 22-27 down to 23-28, preserving their formatting/styles.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly observation.
$ws.Range("A22").Value2 = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value2 = 44748
$ws.Range("E22").Value2 = 7
$ws.Range("F22").Value2 = 100112040
$ws.Range("G22").Value = "Cilantro"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value2 = 150
$ws.Range("K22").Value2 = 8000
$ws.Range("L22").Value2 = 8000
$ws.Range("M22").Value2 = 8000
$ws.Range("N22").Value = "`$/caja 36 atados"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value2 = 222
$ws.Range("Q22").Value2 = 36
$ws.Range("R22").Value = "Hortaliza"

# Append a brand-new row 29 with another new weekly observation.
$ws.Range("A29").Value2 = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D29").Value2 = 44747
$ws.Range("E29").Value2 = 7
$ws.Range("F29").Value2 = 100112040
$ws.Range("G29").Value = "Cilantro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value2 = 200
$ws.Range("K29").Value2 = 9000
$ws.Range("L29").Value2 = 9000
$ws.Range("M29").Value2 = 9000
$ws.Range("N29").Value = "`$/caja 36 atados"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value2 = 250
$ws.Range("Q29").Value2 = 36
$ws.Range("R29").Value = "Hortaliza"
